$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $cellValue)
    # Force the cell to store as literal text (matches source which used
    # inline/shared string cells even for numeric-looking values), then
    # restore the original (unstyled) cell style so no formatting changes.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $cellValue
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '245.28'
Set-TextCell 'G2' '2'
Set-TextCell 'D3' '23.14'
Set-TextCell 'G3' '2'
Set-TextCell 'D4' '5.402'
Set-TextCell 'G4' '2'
Set-TextCell 'D5' '0.06045'
Set-TextCell 'G5' '2'
Set-TextCell 'D6' '3.401'
Set-TextCell 'G6' '2'
Set-TextCell 'D7' '0.8106'
Set-TextCell 'G7' '2'
Set-TextCell 'D8' '0.9334'
Set-TextCell 'G8' '2'
Set-TextCell 'D9' '0.1428'
Set-TextCell 'G9' '2'
Set-TextCell 'D10' '0.07462'
Set-TextCell 'G10' '2'
Set-TextCell 'D11' '0.03345'
Set-TextCell 'G11' '2'
Set-TextCell 'D12' '0.03064'
Set-TextCell 'G12' '2'
Set-TextCell 'D13' '4.010'
Set-TextCell 'G13' '2'
Set-TextCell 'D14' '0.09374'
Set-TextCell 'G14' '2'
Set-TextCell 'D15' '0.001590'
Set-TextCell 'G15' '2'
Set-TextCell 'D16' '0.04832'
Set-TextCell 'G16' '2'
Set-TextCell 'G17' '2'
Set-TextCell 'D18' '0.005426'
Set-TextCell 'G18' '2'
Set-TextCell 'G19' '2'
Set-TextCell 'D20' '0.0009840'
Set-TextCell 'G20' '2'
Set-TextCell 'D21' '0.00008704'
Set-TextCell 'G21' '2'
Set-TextCell 'D22' '3.675'
Set-TextCell 'G22' '2'
Set-TextCell 'D23' '6.430'
Set-TextCell 'G23' '2'
Set-TextCell 'D24' '2.188'
Set-TextCell 'G24' '2'
Set-TextCell 'G25' '2'
Set-TextCell 'D26' '0.1295'
Set-TextCell 'G26' '2'
Set-TextCell 'G27' '2'
Set-TextCell 'G28' '2'
Set-TextCell 'G29' '2'
Set-TextCell 'G30' '2'
Set-TextCell 'G31' '2'
Set-TextCell 'G32' '2'
Set-TextCell 'G33' '2'
Set-TextCell 'G34' '2'
Set-TextCell 'G35' '2'
Set-TextCell 'G36' '2'
Set-TextCell 'G37' '2'
Set-TextCell 'G38' '2'
Set-TextCell 'G39' '2'
Set-TextCell 'D40' '0.03981'
Set-TextCell 'G40' '2'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell 'D41' '0.006403'
$ws.Range('E41').Value = '40KickTokenKICK'
Set-TextCell 'G41' '2'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell 'D42' '0.1075'
$ws.Range('E42').Value = '41BKEXTokenBKK'
Set-TextCell 'G42' '2'
$ws.Range('B43').Value = 'CEJI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell 'D43' '0.002901'
$ws.Range('E43').Value = '42CEJICEJI'
Set-TextCell 'G43' '2'
Set-TextCell 'D44' '0.006347'
Set-TextCell 'G44' '2'
Set-TextCell 'D45' '0.00005257'
Set-TextCell 'G45' '2'
Set-TextCell 'G46' '2'
Set-TextCell 'D47' '0.0005802'
Set-TextCell 'G47' '2'
Set-TextCell 'D48' '0.8869'
Set-TextCell 'G48' '2'
Set-TextCell 'D49' '0.002231'
Set-TextCell 'G49' '2'
Set-TextCell 'G50' '2'
Set-TextCell 'G51' '2'
